$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# New "Lot No." column goes in before the old column E (Result).
$ws.Columns("E:E").Insert()
# New "Remarks" column goes in before the old column H (Created By), which is
# now column I after the first insertion above shifted everything right by one.
$ws.Columns("I:I").Insert()

# --- Column widths ----------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 17.85546875
$ws.Columns("E:E").ColumnWidth = 9.5703125
$ws.Columns("I:I").ColumnWidth = 21.140625

# --- New header cell text ----------------------------------------------------
$ws.Range("E8").Value = "Lot No."
$ws.Range("I8").Value = "Remarks"

# --- Data cell E9:E11 picks up the same format as B9:D11 (done automatically by
# the column insert "format from left" behaviour, nothing further to do) ------

# --- Selection ---------------------------------------------------------------
$ws.Range("B9").Select()

Write-Output "done"
